$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 54549
$ws.Range("E2").Value = 3703
$ws.Range("F2").Value = 3703
$ws.Range("G2").Value = 3709
$ws.Range("H2").Value = 2904
$ws.Range("I2").Value = 2754
$ws.Range("J2").Value = 151
$ws.Range("K2").Value = 32400
$ws.Range("L2").Value = 15700
$ws.Range("M2").Value = 16700
$ws.Range("N2").Value = 16185
$ws.Range("O2").Value = 515
$ws.Range("P2").Value = 534
$ws.Range("Q2").Value = 3468
$ws.Range("R2").Value = -2820
$ws.Range("S2").Value = -440
$ws.Range("T2").Value = 2694
$ws.Range("U2").Value = 775
$ws.Range("V2").Value = 4014
$ws.Range("W2").Value = 6.79
$ws.Range("X2").Value = 5.33
$ws.Range("Y2").Value = 17.7
$ws.Range("Z2").Value = 9.380000000000001
$ws.Range("AA2").Value = 94.01000000000001
$ws.Range("AB2").Value = 3136.05
$ws.Range("AC2").Value = 516
$ws.Range("AD2").Value = 18.76
$ws.Range("AE2").Value = 3032
$ws.Range("AF2").Value = 3.19
$ws.Range("AG2").Value = 194
$ws.Range("AH2").Value = 2
$ws.Range("AI2").Value = 37.6
$ws.Range("AJ2").Value = 533800000

# Row 3
$ws.Range("D3").Value = 55581
$ws.Range("E3").Value = 3596
$ws.Range("F3").Value = 3596
$ws.Range("G3").Value = 3457
$ws.Range("H3").Value = 2434
$ws.Range("I3").Value = 2305
$ws.Range("J3").Value = 129
$ws.Range("K3").Value = 34594
$ws.Range("L3").Value = 16599
$ws.Range("M3").Value = 17994
$ws.Range("N3").Value = 17447
$ws.Range("O3").Value = 547
$ws.Range("P3").Value = 534
$ws.Range("Q3").Value = 4475
$ws.Range("R3").Value = -2559
$ws.Range("S3").Value = -1118
$ws.Range("T3").Value = 2269
$ws.Range("U3").Value = 2206
$ws.Range("V3").Value = 4011
$ws.Range("W3").Value = 6.47
$ws.Range("X3").Value = 4.38
$ws.Range("Y3").Value = 13.71
$ws.Range("Z3").Value = 7.27
$ws.Range("AA3").Value = 92.25
$ws.Range("AB3").Value = 3384.36
$ws.Range("AC3").Value = 432
$ws.Range("AD3").Value = 24.04
$ws.Range("AE3").Value = 3268
$ws.Range("AF3").Value = 3.18
$ws.Range("AG3").Value = 194
$ws.Range("AH3").Value = 1.87
$ws.Range("AI3").Value = 44.93
$ws.Range("AJ3").Value = 533800000

# Row 4
$ws.Range("D4").Value = 57037
$ws.Range("E4").Value = 4225
$ws.Range("F4").Value = 4225
$ws.Range("G4").Value = 4157
$ws.Range("H4").Value = 3038
$ws.Range("I4").Value = 2921
$ws.Range("J4").Value = 116
$ws.Range("K4").Value = 38597
$ws.Range("L4").Value = 19748
$ws.Range("M4").Value = 18849
$ws.Range("N4").Value = 18276
$ws.Range("O4").Value = 573
$ws.Range("P4").Value = 534
$ws.Range("Q4").Value = 3892
$ws.Range("R4").Value = -4519
$ws.Range("S4").Value = 483
$ws.Range("T4").Value = 2497
$ws.Range("U4").Value = 1394
$ws.Range("V4").Value = 6462
$ws.Range("W4").Value = 7.41
$ws.Range("X4").Value = 5.33
$ws.Range("Y4").Value = 16.36
$ws.Range("Z4").Value = 8.300000000000001
$ws.Range("AA4").Value = 104.77
$ws.Range("AB4").Value = 3584.71
$ws.Range("AC4").Value = 547
$ws.Range("AD4").Value = 18.82
$ws.Range("AE4").Value = 3424
$ws.Range("AF4").Value = 3.01
$ws.Range("AG4").Value = 225
$ws.Range("AH4").Value = 2.18
$ws.Range("AI4").Value = 41.11
$ws.Range("AJ4").Value = 533800000

# Row 5
$ws.Range("D5").Value = 55857
$ws.Range("E5").Value = 4684
$ws.Range("F5").Value = 4684
$ws.Range("G5").Value = 4215
$ws.Range("H5").Value = 2984
$ws.Range("I5").Value = 2886
$ws.Range("J5").Value = 98
$ws.Range("K5").Value = 41196
$ws.Range("L5").Value = 20881
$ws.Range("M5").Value = 20314
$ws.Range("N5").Value = 19787
$ws.Range("O5").Value = 528
$ws.Range("P5").Value = 534
$ws.Range("Q5").Value = 5667
$ws.Range("R5").Value = -3465
$ws.Range("S5").Value = -562
$ws.Range("T5").Value = 3221
$ws.Range("U5").Value = 2445
$ws.Range("V5").Value = 7569
$ws.Range("W5").Value = 8.390000000000001
$ws.Range("X5").Value = 5.34
$ws.Range("Y5").Value = 15.16
$ws.Range("Z5").Value = 7.48
$ws.Range("AA5").Value = 102.79
$ws.Range("AB5").Value = 3843.05
$ws.Range("AC5").Value = 541
$ws.Range("AD5").Value = 25.71
$ws.Range("AE5").Value = 3707
$ws.Range("AF5").Value = 3.75
$ws.Range("AG5").Value = 305
$ws.Range("AH5").Value = 2.19
$ws.Range("AI5").Value = 56.42
$ws.Range("AJ5").Value = 533800000

# Row 6
$ws.Range("D6").Value = 59376
$ws.Range("E6").Value = 4338
$ws.Range("F6").Value = 4338
$ws.Range("G6").Value = 3796
$ws.Range("H6").Value = 2837
$ws.Range("I6").Value = 2776
$ws.Range("K6").Value = 53558
$ws.Range("L6").Value = 32241
$ws.Range("M6").Value = 21317
$ws.Range("N6").Value = 20367
$ws.Range("P6").Value = 534
$ws.Range("Q6").Value = 4713
$ws.Range("R6").Value = -5767
$ws.Range("S6").Value = 5637
$ws.Range("T6").Value = 4197
$ws.Range("U6").Value = 515
$ws.Range("V6").Value = 16720
$ws.Range("W6").Value = 7.31
$ws.Range("X6").Value = 4.78
$ws.Range("Y6").Value = 13.83
$ws.Range("Z6").Value = 5.99
$ws.Range("AA6").Value = 151.25
$ws.Range("AB6").Value = 4042.99
$ws.Range("AC6").Value = 520
$ws.Range("AD6").Value = 20.76
$ws.Range("AE6").Value = 3816
$ws.Range("AF6").Value = 2.83
$ws.Range("AG6").Value = 320
$ws.Range("AH6").Value = 2.96
$ws.Range("AI6").Value = 61.52
$ws.Range("AJ6").Value = 533800000

# Row 7
$ws.Range("D7").Value = 71765
$ws.Range("E7").Value = 4698
$ws.Range("G7").Value = 4045
$ws.Range("H7").Value = 3085
$ws.Range("I7").Value = 3067
$ws.Range("K7").Value = 70287
$ws.Range("L7").Value = 46994
$ws.Range("M7").Value = 23293
$ws.Range("N7").Value = 22164
$ws.Range("P7").Value = 532
$ws.Range("Q7").Value = 5901
$ws.Range("R7").Value = -16450
$ws.Range("S7").Value = 8897
$ws.Range("T7").Value = 4903
$ws.Range("U7").Value = 214
$ws.Range("W7").Value = 6.55
$ws.Range("X7").Value = 4.3
$ws.Range("Y7").Value = 14.42
$ws.Range("Z7").Value = 4.98
$ws.Range("AA7").Value = 201.75
$ws.Range("AC7").Value = 574
$ws.Range("AD7").Value = 18.28
$ws.Range("AE7").Value = 4152
$ws.Range("AF7").Value = 2.53
$ws.Range("AG7").Value = 323
$ws.Range("AH7").Value = 3.08
$ws.Range("AI7").Value = 56.24

# Row 8
$ws.Range("D8").Value = 78525
$ws.Range("E8").Value = 5339
$ws.Range("G8").Value = 4733
$ws.Range("H8").Value = 3578
$ws.Range("I8").Value = 3548
$ws.Range("K8").Value = 72364
$ws.Range("L8").Value = 47442
$ws.Range("M8").Value = 24922
$ws.Range("N8").Value = 23743
$ws.Range("P8").Value = 532
$ws.Range("Q8").Value = 7852
$ws.Range("R8").Value = -5346
$ws.Range("S8").Value = -2542
$ws.Range("T8").Value = 3991
$ws.Range("U8").Value = 4014
$ws.Range("W8").Value = 6.8
$ws.Range("X8").Value = 4.56
$ws.Range("Y8").Value = 15.46
$ws.Range("Z8").Value = 5.02
$ws.Range("AA8").Value = 190.36
$ws.Range("AC8").Value = 665
$ws.Range("AD8").Value = 15.8
$ws.Range("AE8").Value = 4448
$ws.Range("AF8").Value = 2.36
$ws.Range("AG8").Value = 345
$ws.Range("AH8").Value = 3.28
$ws.Range("AI8").Value = 51.86

# Row 9
$ws.Range("D9").Value = 82521
$ws.Range("E9").Value = 5941
$ws.Range("G9").Value = 5468
$ws.Range("H9").Value = 4113
$ws.Range("I9").Value = 4056
$ws.Range("K9").Value = 75293
$ws.Range("L9").Value = 48266
$ws.Range("M9").Value = 27027
$ws.Range("N9").Value = 25820
$ws.Range("P9").Value = 532
$ws.Range("Q9").Value = 8494
$ws.Range("R9").Value = -5520
$ws.Range("S9").Value = -2300
$ws.Range("T9").Value = 4105
$ws.Range("U9").Value = 4707
$ws.Range("W9").Value = 7.2
$ws.Range("X9").Value = 4.98
$ws.Range("Y9").Value = 16.37
$ws.Range("Z9").Value = 5.57
$ws.Range("AA9").Value = 178.58
$ws.Range("AC9").Value = 760
$ws.Range("AD9").Value = 13.82
$ws.Range("AE9").Value = 4837
$ws.Range("AF9").Value = 2.17
$ws.Range("AG9").Value = 352
$ws.Range("AH9").Value = 3.36
$ws.Range("AI9").Value = 46.39
